$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.512.17"
$ws.Range("E2").Value = "  +0.49%  "

$ws.Range("D3").Value = "2.439.06"
$ws.Range("E3").Value = "  +1.07%  "

$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "566.13"
$ws.Range("E5").Value = "  +0.70%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.38"
$ws.Range("E6").Value = "  +1.71%  "

$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.533"
$ws.Range("E8").Value = "  +0.44%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.111"
$ws.Range("E9").Value = "  +1.91%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.154"
$ws.Range("E10").Value = "  +0.27%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.25"
$ws.Range("E11").Value = "  -1.30%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.352"
$ws.Range("E12").Value = "  -0.15%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "26.84"
$ws.Range("E13").Value = "  +4.68%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000183"
$ws.Range("E14").Value = "  +4.35%  "

$ws.Range("D15").Value = "2.827.67"
$ws.Range("E15").Value = "  -0.87%  "

$ws.Range("D16").Value = "62.398.87"
$ws.Range("E16").Value = "  +0.55%  "

$ws.Range("D17").Value = "2.440.08"
$ws.Range("E17").Value = "  +1.16%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.25"
$ws.Range("E18").Value = "  -0.47%  "

$ws.Range("E19").Value = "  +1.44%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "324.24"
$ws.Range("E20").Value = "  +0.14%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.17"
$ws.Range("E21").Value = "  -0.22%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.00"
$ws.Range("E22").Value = "  -0.03%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "67.38"
$ws.Range("E23").Value = "  +2.58%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.75"
$ws.Range("E24").Value = "  +2.30%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.81"
$ws.Range("E25").Value = "  -2.37%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "562.08"
$ws.Range("E26").Value = "  -2.50%  "

$ws.Range("D27").Value = "0.0₃0983"
$ws.Range("E27").Value = "  +3.47%  "

$ws.Range("D28").Value = "2.560.11"
$ws.Range("E28").Value = "  +1.14%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  -0.17%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.35"
$ws.Range("E30").Value = "  +1.51%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.45"
$ws.Range("E31").Value = "  +1.11%  "

$ws.Range("E32").Value = "  -0.71%  "

$ws.Range("E33").Value = "  +0.51%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.54"
$ws.Range("E34").Value = "  +0.76%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.88"
$ws.Range("E35").Value = "  +3.47%  "

$ws.Range("E36").Value = "  -0.04%  "

$ws.Range("E37").Value = "  +0.28%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.51"
$ws.Range("E38").Value = "  -0.82%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.76"
$ws.Range("E39").Value = "  +0.45%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "150.53"
$ws.Range("E40").Value = "  -1.09%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.81"
$ws.Range("E41").Value = "  +0.36%  "

$ws.Range("E42").Value = "  +0.75%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.40"
$ws.Range("E43").Value = "  +4.36%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "148.98"
$ws.Range("E44").Value = "  +0.21%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.68"
$ws.Range("E45").Value = "  +0.92%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0536"
$ws.Range("E46").Value = "  +0.27%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "20.41"
$ws.Range("E47").Value = "  +1.84%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.601"
$ws.Range("E48").Value = "  +1.17%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0928"
$ws.Range("E49").Value = "  +1.01%  "

$ws.Range("E50").Value = "  +1.59%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "11.60"
$ws.Range("E51").Value = "  +0.49%  "
